# Workbook: Annotations/Old/1984.xlsx
# Changes made by the author:
#  - Active/selected sheet moved from "Formatted" (3rd tab) to "1948" (1st tab)
#  - On sheet "1948": columns A and B un-hidden and resized
#  - The RAND() driven column A values recompute naturally on save/recalc

$wb = $excel.ActiveWorkbook

$wsOld    = $wb.Worksheets.Item("1948")
$wsChars  = $wb.Worksheets.Item("Characters")
$wsFmt    = $wb.Worksheets.Item("Formatted")

# Unhide & resize column A (was width 11.5546875, hidden) and column B
# (was width 0, hidden) on the "1948" sheet.
$wsOld.Columns.Item(1).Hidden = $false
$wsOld.Columns.Item(1).ColumnWidth = 9
$wsOld.Columns.Item(2).Hidden = $false
$wsOld.Columns.Item(2).ColumnWidth = 13.5

# Make "1948" the active/selected sheet (it was "Formatted" before).
$wsOld.Activate()
